$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("final_figures/2005/F1_P24_Baker & Nelson_2005_ASQ_Creating Something from Nothing.png","Yes","Yes","No","No","No","No","No"),
    @("final_figures/2005/F1_P4_Knight & Pye_2005_HR_Network Learning - an empirically derived model of learning by groups of organizations.png","Yes","Yes","No","No","No","Yes","Yes"),
    @("final_figures/2005/F2_P14_Knight & Pye_2005_HR_Network Learning - an empirically derived model of learning by groups of organizations.png","Yes","No","No","Yes","No","Yes","No"),
    @("final_figures/2005/F1_P15_Wong_2005_OrgSci_Organizational Learning via Expatriate Managers.png","Yes","Yes","No","Yes","No","Yes","No"),
    @("final_figures/2005/F4_P15_Yu, Engleman & Van de Van_2005_OrgStudies_The Integration Journey.png","Yes","Yes","Yes","Yes","No","Yes","No"),
    @("final_figures/2005/F1_P19_Jack_2005_J Management Studies _The Role  Use and Activation of Strong and Weak Network Ties  A Qualitative Analysis.png","Yes","No","No","Yes","No","No","Yes"),
    @("final_figures/2005/F1_P16_Schweizer_2005_AMJ_Organizational Integration of Acquired Biotech Companies.png","Yes","Yes","No","No","No","Yes","No"),
    @("final_figures/2005/F1_P14_Ulrich_2005_J Management StudiesContinuity and Change in Mergers and Acquisitions  A Social Identity Case Study of.png","Yes","Yes","No","No","No","Yes","Yes")
)

$rowIndex = 2
foreach ($row in $data) {
    $colIndex = 1
    foreach ($val in $row) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $val
        $colIndex++
    }
    $rowIndex++
}
